$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (P1, Q1) continuing the numeric sequence,
# copying the bold/bordered header formatting from the adjacent O1 cell.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Update the data rows: swap the I/K and M/O value pairs, and populate the
# two new trailing columns (P, Q) with the same "2" filler value used
# elsewhere in each row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
